$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the newly-used "duplicate_image_filename" column (E) with "NA"
# for every data row (practice rows 2-5 and stimuli rows 6-21).
$ws.Range("E2:E21").Value = "NA"
